$wb = $excel.ActiveWorkbook

$wsBacklog = $wb.Worksheets.Item("Product backlog")
$wsChangeLog = $wb.Worksheets.Item("ChangeLog")

# --- Sheet "ChangeLog": fix typo in C6 (set first so this shared string is added first) ---
$wsChangeLog.Range("C6").Value = "verbeteringen product backlog"

# --- Sheet "Product backlog": reword acceptance criteria in A2 ---
$wsBacklog.Range("A2").Value = "1. Het systeem heeft een interface met sensordata waarin de data overzichtelijk voor weergeven. Bijv. een grafiek.  "

# --- Update selection / view state (select ChangeLog first, then Product backlog last so it stays the active/tabSelected sheet) ---
$wsChangeLog.Range("C6").Select()
$wsBacklog.Range("D3").Select()
